# Update ticket/interest counts ("F" column) for two events in two sheets:
#   展览   (Exhibitions)  - sheet with these rows
#   全部类型 (All types)   - combined sheet with these rows

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 16509
$wsExpo.Range("F4").Value = 742
$wsExpo.Range("F6").Value = 710
$wsExpo.Range("F7").Value = 1758

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 16509
$wsAll.Range("F4").Value = 742
$wsAll.Range("F8").Value = 710
$wsAll.Range("F9").Value = 1758
